$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.640.66"
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = "'2.022.88"
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'263.60"
$ws.Range('E5').Value = '  +6.43%  '
$ws.Range('E6').Value = '  -1.62%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'56.48"
$ws.Range('E8').Value = '  -6.75%  '
$ws.Range('D9').Value = "'0.387"
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('E10').Value = '  -2.77%  '
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('E12').Value = '  -4.01%  '
$ws.Range('D13').Value = "'2.320.10"
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = "'0.811"
$ws.Range('E14').Value = '  -4.47%  '
$ws.Range('D15').Value = "'20.89"
$ws.Range('E15').Value = '  -8.37%  '
$ws.Range('D16').Value = "'5.28"
$ws.Range('E16').Value = '  -3.16%  '
$ws.Range('D17').Value = "'2.030.56"
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = "'37.462.44"
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').Value = "'69.94"
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = "'0.0₃0847"
$ws.Range('E20').Value = '  -2.21%  '
$ws.Range('D21').Value = "'5.19"
$ws.Range('D22').Value = "'229.08"
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').Value = "'2.71"
$ws.Range('E23').Value = '  +7.91%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = "'2.33"
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('D26').Value = "'164.29"
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').Value = "'9.02"
$ws.Range('E27').Value = '  -4.28%  '
$ws.Range('D28').Value = "'19.76"
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  -11.01%  '
$ws.Range('D30').Value = "'1.34"
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').Value = "'0.0654"
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E33').Value = '  -3.92%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').Value = "'2.39"
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').Value = "'1.82"
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').Value = "'3.37"
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('D39').Value = "'5.24"
$ws.Range('E39').Value = '  -5.02%  '
$ws.Range('D40').Value = "'3.07"
$ws.Range('E40').Value = '  +4.66%  '
$ws.Range('E41').Value = '  +3.61%  '
$ws.Range('D42').Value = "'0.0943"
$ws.Range('E42').Value = '  -3.45%  '
$ws.Range('D43').Value = "'0.0214"
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').Value = "'1.408.31"
$ws.Range('D45').Value = "'90.93"
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = "'15.77"
$ws.Range('E46').Value = '  -5.58%  '
$ws.Range('D47').Value = "'1.03"
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = "'7.10"
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').Value = "'2.88"
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = "'1.98"
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = "'2.211.09"
$ws.Range('E51').Value = '  +0.81%  '
